# Update the "2024" worksheet: a new SMS-parsed transaction entry was
# recorded ahead of the existing September entries, so insert a new row
# at row 29 (pushing the existing rows 29:43 down to 30:44) and populate
# the new row's September Details / September Date cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Rows.Item(29).Insert()

$ws.Range("R29").Value = "login sbi internet personal do not share anyone"
$ws.Range("S29").Value = "2024-09-03 19:17:10"
